# Update LR-pair rows to reflect revised Il1b/Il1r1 Sending-cluster analysis
# (ECs and sCs as sending clusters, each paired against ECs/FAPs/sCs targets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1b"
$ws.Range("C2").Value = "Il1r1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1413.335253
$ws.Range("H2").Value = 4240.005759
$ws.Range("I2").Value = 0.9999668843963775
$ws.Range("J2").Value = 0.9999668843963775
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.175804666666667
$ws.Range("N2").Value = 21.527414
$ws.Range("O2").Value = 0.1618789277039843
$ws.Range("P2").Value = 0.1618789277039842
$ws.Range("Q2").Value = 10141.81770404191
$ws.Range("R2").Value = 91276.35933637722
$ws.Range("S2").Value = 0.1618735669855796
$ws.Range("T2").Value = 0.1618735669855796

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1b"
$ws.Range("C3").Value = "Il1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1413.335253
$ws.Range("H3").Value = 4240.005759
$ws.Range("I3").Value = 0.9999668843963775
$ws.Range("J3").Value = 0.9999668843963775
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.92162533333333
$ws.Range("N3").Value = 74.764876
$ws.Range("O3").Value = 0.5622067730383848
$ws.Range("P3").Value = 0.5622067730383847
$ws.Range("Q3").Value = 35222.61164565787
$ws.Range("R3").Value = 317003.5048109209
$ws.Range("S3").Value = 0.562188155221735
$ws.Range("T3").Value = 0.5621881552217349

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1b"
$ws.Range("C4").Value = "Il1r1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1413.335253
$ws.Range("H4").Value = 4240.005759
$ws.Range("I4").Value = 0.9999668843963775
$ws.Range("J4").Value = 0.9999668843963775
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.23078966666667
$ws.Range("N4").Value = 36.692369
$ws.Range("O4").Value = 0.2759142992576309
$ws.Range("P4").Value = 0.2759142992576308
$ws.Range("Q4").Value = 17286.20620792812
$ws.Range("R4").Value = 155575.8558713531
$ws.Range("S4").Value = 0.2759051621890629
$ws.Range("T4").Value = 0.2759051621890629

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Il1b"
$ws.Range("C5").Value = "Il1r1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04680500000000001
$ws.Range("H5").Value = 0.140415
$ws.Range("I5").Value = [double]"3.311560362258399e-05"
$ws.Range("J5").Value = [double]"3.311560362258399e-05"
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.175804666666667
$ws.Range("N5").Value = 21.527414
$ws.Range("O5").Value = 0.1618789277039843
$ws.Range("P5").Value = 0.1618789277039842
$ws.Range("Q5").Value = 0.3358635374233334
$ws.Range("R5").Value = 3.022771836810001
$ws.Range("S5").Value = [double]"5.360718404694073e-06"
$ws.Range("T5").Value = [double]"5.360718404694072e-06"

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Il1b"
$ws.Range("C6").Value = "Il1r1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04680500000000001
$ws.Range("H6").Value = 0.140415
$ws.Range("I6").Value = [double]"3.311560362258399e-05"
$ws.Range("J6").Value = [double]"3.311560362258399e-05"
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.92162533333333
$ws.Range("N6").Value = 74.764876
$ws.Range("O6").Value = 0.5622067730383848
$ws.Range("P6").Value = 0.5622067730383847
$ws.Range("Q6").Value = 1.166456673726667
$ws.Range("R6").Value = 10.49811006354
$ws.Range("S6").Value = [double]"1.861781664987119e-05"
$ws.Range("T6").Value = [double]"1.861781664987119e-05"

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Il1b"
$ws.Range("C7").Value = "Il1r1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04680500000000001
$ws.Range("H7").Value = 0.140415
$ws.Range("I7").Value = [double]"3.311560362258399e-05"
$ws.Range("J7").Value = [double]"3.311560362258399e-05"
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.23078966666667
$ws.Range("N7").Value = 36.692369
$ws.Range("O7").Value = 0.2759142992576309
$ws.Range("P7").Value = 0.2759142992576308
$ws.Range("Q7").Value = 0.5724621103483334
$ws.Range("R7").Value = 5.152158993135
$ws.Range("S7").Value = [double]"9.137068568018724e-06"
$ws.Range("T7").Value = [double]"9.137068568018722e-06"
